# Power Models and Textures.xlsx - add 4 new texture entries to Table2 on Sheet1
# (the table is sorted A-Z by the "File" column, so after adding the rows we
# re-apply the table's sort to land them in the correct alphabetical slot.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table2")

$newEntries = @(
    @("textures\grass.igb", "XML2 PC", "XML2 GameCube", "XML2 PS2", "XML2 PSP", "XML2 Xbox", "1b. Used in all versions of XML2 but not originally in permanent"),
    @("textures\water_insect.igb", "XML2 PC", "XML2 GameCube", "XML2 PS2", "XML2 PSP", "XML2 Xbox", "1b. Used in all versions of XML2 but not originally in permanent"),
    @("textures\water_sewers6.igb", "XML2 PC", "XML2 GameCube", "XML2 PS2", "XML2 PSP", "XML2 Xbox", "1b. Used in all versions of XML2 but not originally in permanent"),
    @("textures\debris\paper_4.igb", "XML2 PC", "XML2 GameCube", "XML2 PS2", "XML2 PSP", "XML2 Xbox", "1b. Used in all versions of XML2 but not originally in permanent")
)

foreach ($entry in $newEntries) {
    $row = $lo.ListRows.Add()
    $row.Range.Item(1,1).Value = $entry[0]
    $row.Range.Item(1,2).Value = $entry[1]
    $row.Range.Item(1,3).Value = $entry[2]
    $row.Range.Item(1,4).Value = $entry[3]
    $row.Range.Item(1,5).Value = $entry[4]
    $row.Range.Item(1,6).Value = $entry[5]
    $row.Range.Item(1,7).Value = $entry[6]
}

# The table keeps an ascending sort on column A ("File"); re-apply it so the
# newly-appended rows move into their correct alphabetical position.
$lo.Sort.SortFields.Clear()
[void]$lo.Sort.SortFields.Add($lo.ListColumns.Item(1).Range)
$lo.Sort.Header = 1
$lo.Sort.Apply()

# Leave the selection where the author's cursor ended up after entering the
# last new row (one past the bottom-right of the grown table).
[void]$ws.Range("G144").Select()
